$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4874
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 4874
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 14622
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -15184

$ws.Range("H55").Value = 849.8333
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 919.8
$ws.Range("K55").Value = 500
$ws.Range("L55").Value = 919.8
$ws.Range("M55").Value = -286
$ws.Range("N55").Value = -1347.8

$ws.Range("H76").Value = 7958.3335
$ws.Range("I76").Value = 7750
$ws.Range("J76").Value = 8000
$ws.Range("K76").Value = 7750
$ws.Range("L76").Value = 8000
$ws.Range("M76").Value = -7435
$ws.Range("N76").Value = -8630

$ws.Range("H79").Value = 7958.3335
$ws.Range("I79").Value = 7750
$ws.Range("J79").Value = 8000
$ws.Range("K79").Value = 7750
$ws.Range("L79").Value = 8000
$ws.Range("M79").Value = -6658
$ws.Range("N79").Value = -10184

$ws.Range("H86").Value = 4844.1113
$ws.Range("I86").Value = 2532.6667
$ws.Range("J86").Value = 5999.8335
$ws.Range("K86").Value = 2532.6667
$ws.Range("L86").Value = 5999.8335
$ws.Range("M86").Value = -1409.6667
$ws.Range("N86").Value = -8245.833500000001

$ws.Range("H89").Value = 4844.1113
$ws.Range("I89").Value = 2532.6667
$ws.Range("J89").Value = 5999.8335
$ws.Range("K89").Value = 12663.3335
$ws.Range("L89").Value = 29999.1675
$ws.Range("M89").Value = -7047.333500000001
$ws.Range("N89").Value = -41231.1675

$ws.Range("H93").Value = 34500
$ws.Range("J93").Value = 34500
$ws.Range("L93").Value = 34500
$ws.Range("N93").Value = -39492

$ws.Range("H98").Value = 1249.3077
$ws.Range("I98").Value = 1263.2727
$ws.Range("K98").Value = 1263.2727
$ws.Range("M98").Value = 234.7273

$ws.Range("H106").Value = 24071.715
$ws.Range("I106").Value = 34668.777
$ws.Range("J106").Value = 4997
$ws.Range("K106").Value = 34668.777
$ws.Range("L106").Value = 4997
$ws.Range("M106").Value = -34037.777
$ws.Range("N106").Value = -6259

$ws.Range("H122").Value = 1249.3077
$ws.Range("I122").Value = 1263.2727
$ws.Range("K122").Value = 3789.8181
$ws.Range("M122").Value = -1339.8181

$ws.Range("H132").Value = 1941.6
$ws.Range("I132").Value = 1899.7646
$ws.Range("J132").Value = 2178.6667
$ws.Range("K132").Value = 5699.293799999999
$ws.Range("L132").Value = 6536.000100000001
$ws.Range("M132").Value = -3169.293799999999
$ws.Range("N132").Value = -11596.0001

$ws.Range("H135").Value = 1044.2727
$ws.Range("I135").Value = 811.1875
$ws.Range("K135").Value = 7300.6875
$ws.Range("M135").Value = -4765.6875

$ws.Range("H137").Value = 5616.8887
$ws.Range("I137").Value = 2250.3333
$ws.Range("J137").Value = 8983.444
$ws.Range("K137").Value = 6750.999899999999
$ws.Range("L137").Value = 26950.332
$ws.Range("M137").Value = -4200.999899999999
$ws.Range("N137").Value = -32050.332

$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 4998
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4998
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 4998
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -5366

$ws.Range("H32").Value = 15360.947
$ws.Range("I32").Value = 6045.114
$ws.Range("J32").Value = 28170.219
$ws.Range("K32").Value = 6045.114
$ws.Range("L32").Value = 28170.219
$ws.Range("M32").Value = -5758.114
$ws.Range("N32").Value = -28744.219

$ws.Range("H122").Value = 336428.84
$ws.Range("I122").Value = 502394.4
$ws.Range("K122").Value = 1507183.2
$ws.Range("M122").Value = -1504733.2

$ws.Range("H132").Value = 3856.8333
$ws.Range("I132").Value = 1598.2354
$ws.Range("K132").Value = 4794.706200000001
$ws.Range("M132").Value = -2264.706200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2031.9231
$ws.Range("J99").Value = 2113
$ws.Range("L99").Value = 2113
$ws.Range("N99").Value = -5109

$ws.Range("H105").Value = 4314.6895
$ws.Range("I105").Value = 3627.7646
$ws.Range("K105").Value = 3627.7646
$ws.Range("M105").Value = -1880.7646

$ws.Range("H107").Value = 2724.85
$ws.Range("J107").Value = 4570.857
$ws.Range("L107").Value = 4570.857
$ws.Range("N107").Value = -8410.857

$ws.Range("H134").Value = 2924.1304
$ws.Range("I134").Value = 2119.1052
$ws.Range("K134").Value = 6357.3156
$ws.Range("M134").Value = -3822.3156

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4276.0835
$ws.Range("I31").Value = 3013.5454
$ws.Range("K31").Value = 3013.5454
$ws.Range("M31").Value = -2718.5454

$ws.Range("H34").Value = 4276.0835
$ws.Range("I34").Value = 3013.5454
$ws.Range("K34").Value = 3013.5454
$ws.Range("M34").Value = -2811.5454

$ws.Range("H122").Value = 915.5
$ws.Range("I122").Value = 923.75
$ws.Range("K122").Value = 2771.25
$ws.Range("M122").Value = -321.25

$ws.Range("H132").Value = 1931.2188
$ws.Range("I132").Value = 1867.8889
$ws.Range("J132").Value = 2273.2
$ws.Range("K132").Value = 5603.6667
$ws.Range("L132").Value = 6819.599999999999
$ws.Range("M132").Value = -3073.6667
$ws.Range("N132").Value = -11879.6

$ws.Range("H134").Value = 3669.8125
$ws.Range("I134").Value = 2434.2727
$ws.Range("K134").Value = 7302.8181
$ws.Range("M134").Value = -4767.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2535.4443
$ws.Range("I14").Value = 2535.4443
$ws.Range("K14").Value = 7606.3329
$ws.Range("M14").Value = -7433.3329

$ws.Range("H23").Value = 172.16667
$ws.Range("I23").Value = 172.16667
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 516.50001
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -281.50001
$ws.Range("N23").ClearContents()

$ws.Range("H34").Value = 2300
$ws.Range("J34").Value = 2665
$ws.Range("L34").Value = 7995
$ws.Range("N34").Value = -8163

$ws.Range("H39").Value = 3855.5557
$ws.Range("J39").Value = 4950
$ws.Range("L39").Value = 14850
$ws.Range("N39").Value = -15438

$ws.Range("H55").Value = 1470
$ws.Range("J55").Value = 1587.5
$ws.Range("L55").Value = 4762.5
$ws.Range("N55").Value = -5116.5

$ws.Range("H56").Value = 11906.667
$ws.Range("I56").Value = 11906.667
$ws.Range("K56").Value = 11906.667
$ws.Range("M56").Value = -11376.667

$ws.Range("H131").Value = 1292.7778
$ws.Range("I131").Value = 699.5
$ws.Range("J131").Value = 1462.2858
$ws.Range("K131").Value = 2098.5
$ws.Range("L131").Value = 4386.857400000001
$ws.Range("M131").Value = 2941.5
$ws.Range("N131").Value = -14466.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6234.6924
$ws.Range("I70").Value = 4411.2
$ws.Range("K70").Value = 4411.2
$ws.Range("M70").Value = -4141.2

$ws.Range("H73").Value = 6234.6924
$ws.Range("I73").Value = 4411.2
$ws.Range("K73").Value = 4411.2
$ws.Range("M73").Value = -3475.2

$ws.Range("H102").Value = 2175.6
$ws.Range("I102").Value = 1482.7778
$ws.Range("K102").Value = 1482.7778
$ws.Range("M102").Value = 139.2221999999999

$ws.Range("H126").Value = 7497.5
$ws.Range("J126").Value = 7497.5
$ws.Range("L126").Value = 22492.5
$ws.Range("N126").Value = -27432.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3250.1538
$ws.Range("I61").Value = 2841.4546
$ws.Range("K61").Value = 2841.4546
$ws.Range("M61").Value = -2639.4546

$ws.Range("H113").Value = 3250.1538
$ws.Range("I113").Value = 2841.4546
$ws.Range("K113").Value = 2841.4546
$ws.Range("M113").Value = -671.4546

$ws.Range("H132").Value = 4546.5625
$ws.Range("I132").Value = 3027.2222
$ws.Range("K132").Value = 9081.6666
$ws.Range("M132").Value = -6551.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 2583.3333
$ws.Range("I6").Value = 2375
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 2375
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -2260
$ws.Range("N6").Value = -3230

$ws.Range("H113").Value = 1995.6
$ws.Range("I113").Value = 1175.6666
$ws.Range("K113").Value = 3526.9998
$ws.Range("M113").Value = -1356.9998

$ws.Range("H132").Value = 2720.9
$ws.Range("I132").Value = 2245.25
$ws.Range("J132").Value = 3434.375
$ws.Range("K132").Value = 6735.75
$ws.Range("L132").Value = 10303.125
$ws.Range("M132").Value = -4205.75
$ws.Range("N132").Value = -15363.125

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
